$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.597.46'
$ws.Range("E2").Value = '  -1.06%  '

$ws.Range("D3").Value = '3.780.54'
$ws.Range("E3").Value = '  +1.03%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.77'
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.64'
$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("D7").Value = '3.767.03'
$ws.Range("E7").Value = '  +0.76%  '

$ws.Range("E9").Value = '  +0.40%  '

$ws.Range("E10").Value = '  -0.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.30'
$ws.Range("E11").Value = '  -2.08%  '

$ws.Range("E12").Value = '  -0.11%  '

$ws.Range("E13").Value = '  -2.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.94'
$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("D15").Value = '4.415.90'

$ws.Range("D16").Value = '3.790.21'
$ws.Range("E16").Value = '  +1.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.50'
$ws.Range("E17").Value = '  +3.54%  '

$ws.Range("D18").Value = '67.599.65'
$ws.Range("E18").Value = '  -1.13%  '

$ws.Range("E19").Value = '  +0.22%  '

$ws.Range("E20").Value = '  -0.02%  '

$ws.Range("E21").Value = '  -5.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.16'
$ws.Range("E22").Value = '  -1.37%  '

$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("E24").Value = '  +5.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.37'
$ws.Range("E25").Value = '  -0.70%  '

$ws.Range("E26").Value = '  +0.96%  '

$ws.Range("E27").Value = '  -2.85%  '

$ws.Range("E28").Value = '  +0.11%  '

$ws.Range("E29").Value = '  -0.63%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  +3.12%  '

$ws.Range("E32").Value = '  -1.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.54'
$ws.Range("E33").Value = '  -0.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.37%  '

$ws.Range("E35").Value = '  -1.05%  '

$ws.Range("E36").Value = '  -0.32%  '

$ws.Range("E37").Value = '  -0.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.137'
$ws.Range("E38").Value = '  -0.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.993'
$ws.Range("E39").Value = '  -0.41%  '

$ws.Range("E40").Value = '  -0.37%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '45.43'
$ws.Range("E43").Value = '  +3.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.11'
$ws.Range("E44").Value = '  +3.42%  '

$ws.Range("E46").Value = '  +3.75%  '

$ws.Range("E47").Value = '  -1.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '392.98'
$ws.Range("E48").Value = '  +0.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.57'
$ws.Range("E49").Value = '  +5.72%  '

$ws.Range("E50").Value = '  -5.31%  '

$ws.Range("D51").Value = '2.716.23'
$ws.Range("E51").Value = '  -1.08%  '
